# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric to Excel's auto-detection but must
# remain plain text (matching the workbook's original inlineStr cells).
# Force text format before writing, then restore the default 'Normal' style
# so no stray number-format style is left behind on the cell.
$textForceCells = @('D5', 'D6', 'D8', 'D9', 'D10', 'D11', 'D13', 'D14', 'D16', 'D17', 'D19', 'D21', 'D22', 'D23', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D36', 'D37', 'D38', 'D40', 'D41', 'D42', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell (row order matches the sheet).
$ws.Range('D2').Value = '28.508.33'
$ws.Range('E2').Value = '  -2.46%  '
$ws.Range('D3').Value = '1.792.44'
$ws.Range('E3').Value = '  -1.99%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '231.41'
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('D6').Value = '0.5889'
$ws.Range('E6').Value = '  -1.47%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '0.2768'
$ws.Range('D9').Value = '23.35'
$ws.Range('E9').Value = '  +0.54%  '
$ws.Range('D10').Value = '0.06740'
$ws.Range('E10').Value = '  -3.00%  '
$ws.Range('D11').Value = '0.07569'
$ws.Range('E11').Value = '  -0.60%  '
$ws.Range('D12').Value = '1.795.62'
$ws.Range('E12').Value = '  -1.98%  '
$ws.Range('D13').Value = '4.790'
$ws.Range('E13').Value = '  +0.88%  '
$ws.Range('D14').Value = '0.6114'
$ws.Range('E14').Value = '  -2.14%  '
$ws.Range('D15').Value = '2.035.45'
$ws.Range('E15').Value = '  -1.99%  '
$ws.Range('D16').Value = '75.56'
$ws.Range('E16').Value = '  -3.53%  '
$ws.Range('D17').Value = '0.000008835'
$ws.Range('E17').Value = '  -8.78%  '
$ws.Range('D18').Value = '28.500.38'
$ws.Range('E18').Value = '  -1.52%  '
$ws.Range('D19').Value = '5.435'
$ws.Range('E19').Value = '  -4.87%  '
$ws.Range('E20').Value = '  -0.15%  '
$ws.Range('D21').Value = '208.22'
$ws.Range('E21').Value = '  -6.01%  '
$ws.Range('D22').Value = '11.44'
$ws.Range('E22').Value = '  -0.56%  '
$ws.Range('D23').Value = '6.807'
$ws.Range('E23').Value = '  -0.81%  '
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').Value = '152.78'
$ws.Range('E25').Value = '  -1.95%  '
$ws.Range('D26').Value = '7.973'
$ws.Range('E26').Value = '  +0.56%  '
$ws.Range('D27').Value = '0.1260'
$ws.Range('E27').Value = '  -2.09%  '
$ws.Range('D28').Value = '16.40'
$ws.Range('E28').Value = '  -0.51%  '
$ws.Range('D29').Value = '1.416'
$ws.Range('E29').Value = '  -2.23%  '
$ws.Range('D30').Value = '0.06109'
$ws.Range('E30').Value = '  -7.93%  '
$ws.Range('D31').Value = '1.421'
$ws.Range('D32').Value = '3.780'
$ws.Range('E32').Value = '  -1.21%  '
$ws.Range('D33').Value = '3.761'
$ws.Range('E33').Value = '  +0.46%  '
$ws.Range('D34').Value = '1.721'
$ws.Range('E34').Value = '  +0.61%  '
$ws.Range('E35').Value = '  -3.90%  '
$ws.Range('D36').Value = '0.6417'
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').Value = '2.502'
$ws.Range('E37').Value = '  -1.74%  '
$ws.Range('D38').Value = '2.702'
$ws.Range('E38').Value = '  -1.19%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.149.88'
$ws.Range('E39').Value = '  -2.46%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '0.01682'
$ws.Range('E40').Value = '  -2.78%  '
$ws.Range('D41').Value = '6.309'
$ws.Range('E41').Value = '  -3.22%  '
$ws.Range('D42').Value = '0.8729'
$ws.Range('E42').Value = '  -2.96%  '
$ws.Range('E43').Value = '  -0.11%  '
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('D45').Value = '1.944.82'
$ws.Range('E45').Value = '  -1.76%  '
$ws.Range('D46').Value = '60.17'
$ws.Range('E46').Value = '  -2.85%  '
$ws.Range('D47').Value = '0.00000000111'
$ws.Range('E47').Value = '  -2.59%  '
$ws.Range('D48').Value = '1.584'
$ws.Range('E48').Value = '  +0.68%  '
$ws.Range('D49').Value = '8.344'
$ws.Range('E49').Value = '  -0.96%  '
$ws.Range('D50').Value = '0.05443'
$ws.Range('E50').Value = '  -1.39%  '
$ws.Range('D51').Value = '0.4473'
$ws.Range('E51').Value = '  -1.73%  '

# Restore default styling on the cells we temporarily formatted as text,
# so their style index matches the original (unstyled) cells.
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
